$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 157 is the last existing data row (date serial 45713, 2025-02-25).
# Append 20 more daily rows (158-177, through 2025-03-17): same B:J
# values as row 157, with column A incrementing by one day each row.
$lastRow = 157
$newRows = 20
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

for ($i = 1; $i -le $newRows; $i++) {
    $srcRow = $lastRow
    $dstRow = $lastRow + $i

    $srcRange = "A" + $srcRow + ":J" + $srcRow
    $dstRange = "A" + $dstRow + ":J" + $dstRow

    # Copy the whole source row (values + style) to the new row first, so
    # the date cell keeps the same number-format style, then fix up values.
    $ws.Range($srcRange).Copy($ws.Range($dstRange))

    $ws.Range("A" + $dstRow).Value = $ws.Range("A" + $srcRow).Value2 + $i

    foreach ($col in $cols) {
        $ws.Range($col + $dstRow).Value = $ws.Range($col + $srcRow).Value2
    }
}
